$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column S ("phase"), shifting everything to its right one column to the left.
$ws.Columns("S").Delete()

# Clear the now-orphaned trailing data (previously AD3:AF3, now shifted to AC3:AE3)
# that is no longer part of the sheet's used layout.
$ws.Range("AC3:AE3").ClearContents()

# The "Additional Options" merged banner (was X1:AC1, now W1:AB1 after the column
# delete) grows to also cover the former annual/period_start/period_end columns
# (now AC1:AE1), matching the rest of the header row's style.
$ws.Range("X1:AE1").Style = $ws.Range("C1").Style
$ws.Range("W1:AE1").Merge()
$ws.Range("W1").Value = "Additional Options"

# Match the new active selection left by this edit.
$ws.Range("X5").Select() | Out-Null
